$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.769
$ws.Range("D7").Value = -7.842000000000001
$ws.Range("D16").Value = -8.586
$ws.Range("D28").Value = -8.262
$ws.Range("D29").Value = -7.292
$ws.Range("D32").Value = -7.935999999999998
$ws.Range("D40").Value = -8.132
$ws.Range("D52").Value = -7.87
$ws.Range("D57").Value = -8.101000000000001
$ws.Range("D66").Value = -7.436
$ws.Range("D100").Value = -8.279999999999999
